# Fruta / hortaliza, semanal
# The daily rows of data (D, J, K, L, M, O, P columns) got reshuffled /
# re-associated with different dates (i.e. a permutation of the 55 data
# rows, 2..56). Below we snapshot every row's current values for those
# columns, then write them back out according to the new row order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (source row's CURRENT / BEFORE
# values for columns D, J, K, L, M, O, P become the new values for the
# destination row).
$map = @{
    2 = 22
    3 = 53
    4 = 43
    5 = 38
    6 = 55
    7 = 34
    8 = 23
    9 = 3
    10 = 46
    11 = 16
    12 = 30
    13 = 28
    14 = 54
    15 = 8
    16 = 26
    17 = 37
    18 = 48
    19 = 39
    20 = 50
    21 = 44
    22 = 41
    23 = 49
    24 = 51
    25 = 29
    26 = 27
    27 = 18
    28 = 52
    29 = 12
    30 = 9
    31 = 4
    32 = 7
    33 = 21
    34 = 13
    35 = 56
    36 = 31
    37 = 47
    38 = 45
    39 = 2
    40 = 11
    41 = 33
    42 = 25
    43 = 24
    44 = 40
    45 = 36
    46 = 19
    47 = 10
    48 = 14
    49 = 42
    50 = 15
    51 = 32
    52 = 20
    53 = 35
    54 = 17
    55 = 5
    56 = 6
}

# Columns that move together with each row.
$cols = @("D", "J", "K", "L", "M", "O", "P")

# 1) Snapshot the current ("before") values of every relevant cell.
$snapshot = @{}
for ($r = 2; $r -le 56; $r++) {
    foreach ($c in $cols) {
        $addr = "$c$r"
        $snapshot[$addr] = $ws.Range($addr).Value2
    }
}

# 2) Write the snapshotted values back according to the new mapping.
foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    foreach ($c in $cols) {
        $srcAddr = "$c$srcRow"
        $destAddr = "$c$destRow"
        $ws.Range($destAddr).Value = $snapshot[$srcAddr]
    }
}
